# Insert a new weekly price record for "Ciboulette" (Mercado Mayorista Lo
# Valledor de Santiago) right before the current row 357. All existing
# rows from 357 downward are pushed one row lower (xlShiftDown = -4121),
# which grows the sheet from A1:R381 to A1:R382 and reproduces the row
# renumbering seen throughout the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 357:381 down to 358:382, opening up a blank row 357.
$ws.Rows.Item(357).Insert(-4121)

# Populate the newly inserted row 357 with the new weekly record.
$ws.Range("A357").Value = 6
$ws.Range("B357").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C357").Value = "Metropolitana"
$ws.Range("D357").Value = 44610
$ws.Range("E357").Value = 13
$ws.Range("F357").Value = 100112039
$ws.Range("G357").Value = "Ciboulette"
$ws.Range("H357").Value = "Sin especificar"
$ws.Range("I357").Value = "Primera"
$ws.Range("J357").Value = 820
$ws.Range("K357").Value = 900
$ws.Range("L357").Value = 1000
$ws.Range("M357").Value = 943
$ws.Range("N357").Value = "$/docena de atados"
$ws.Range("O357").Value = "Región Metropolitana"
$ws.Range("P357").Value = 314
$ws.Range("Q357").Value = 3
$ws.Range("R357").Value = "Hortaliza"
